# Updated symbol list on Sun Jan  8 19:31:37 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency rows that changed between the previous and latest scrape.
# Values are written with a leading apostrophe so Excel keeps them as text
# (matching the original inlineStr cells) instead of auto-converting the
# numeric-looking strings/percentages into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'265.82"
$ws.Range("E2").Value  = "'1.64%"

$ws.Range("D3").Value  = "'26.70"
$ws.Range("E3").Value  = "'-1.85%"

$ws.Range("D4").Value  = "'4.706"
$ws.Range("E4").Value  = "'0.03%"

$ws.Range("D5").Value  = "'0.06082"
$ws.Range("E5").Value  = "'-1.82%"

$ws.Range("D6").Value  = "'6.740"
$ws.Range("E6").Value  = "'0.40%"

$ws.Range("D7").Value  = "'0.8514"
$ws.Range("E7").Value  = "'0.10%"

$ws.Range("D8").Value  = "'0.9063"
$ws.Range("E8").Value  = "'-0.96%"

$ws.Range("E9").Value  = "'-0.13%"

$ws.Range("D10").Value = "'0.04998"
$ws.Range("E10").Value = "'9.98%"

$ws.Range("D11").Value = "'0.07090"
$ws.Range("E11").Value = "'0.09%"

$ws.Range("D12").Value = "'0.03188"
$ws.Range("E12").Value = "'1.77%"

$ws.Range("D13").Value = "'0.09012"
$ws.Range("E13").Value = "'-0.44%"

$ws.Range("D14").Value = "'0.001537"
$ws.Range("E14").Value = "'0.60%"

$ws.Range("D15").Value = "'0.0006044"
$ws.Range("E15").Value = "'-1.94%"

$ws.Range("D16").Value = "'0.006016"
$ws.Range("E16").Value = "'0.38%"

$ws.Range("D17").Value = "'3.461"
$ws.Range("E17").Value = "'0.06%"

$ws.Range("D18").Value = "'3.170"
$ws.Range("E18").Value = "'0.12%"

$ws.Range("D19").Value = "'2.264"
$ws.Range("E19").Value = "'3.24%"

$ws.Range("D20").Value = "'0.3089"
$ws.Range("E20").Value = "'-0.57%"

$ws.Range("E21").Value = "'-0.65%"

$ws.Range("D22").Value = "'4.075"
$ws.Range("E22").Value = "'-0.53%"

$ws.Range("D23").Value = "'0.04230"
$ws.Range("E23").Value = "'-0.07%"

$ws.Range("D24").Value = "'0.001182"
$ws.Range("E24").Value = "'-2.94%"

$ws.Range("D25").Value = "'0.004135"
$ws.Range("E25").Value = "'8.76%"

$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.03%"

$ws.Range("E27").Value = "'5.05%"

$ws.Range("D40").Value = "'0.03916"
$ws.Range("E40").Value = "'-0.22%"

$ws.Range("E41").Value = "'-0.04%"

$ws.Range("D42").Value = "'0.004178"
$ws.Range("E42").Value = "'1.14%"

$ws.Range("D43").Value = "'0.002111"
$ws.Range("E43").Value = "'-3.33%"

$ws.Range("D44").Value = "'0.01270"
$ws.Range("E44").Value = "'-8.27%"

$ws.Range("D45").Value = "'0.00005102"
$ws.Range("E45").Value = "'-0.95%"

$ws.Range("E46").Value = "'0.03%"

$ws.Range("D48").Value = "'0.1350"
$ws.Range("E48").Value = "'-19.51%"

$ws.Range("E49").Value = "'0.03%"

$ws.Range("E50").Value = "'0.03%"
